$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the SQL queries (TabQuery column B, StatQuery cell C2) ---
# Each query's join conditions are updated from generic ".id" keys to the
# explicit "<table>_id" keys.

function Update-QueryText([string]$text) {
    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    return $text
}

$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $cell.Value2 = Update-QueryText $cell.Value2
}

# --- Update the sheet view: clear the scrolled topLeftCell and move the
#     active selection from B7 to B2 ---
$ws.Activate()
$ws.Range("B2").Select() | Out-Null
